$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 551; this shifts the existing rows 551:641
# down to 552:642 (and grows the sheet's used range accordingly).
$ws.Rows.Item(551).Insert()

# Populate the newly inserted row 551 with the new weekly record.
$ws.Cells.Item(551, 1).Value = 9
$ws.Cells.Item(551, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(551, 3).Value = "Metropolitana"
$ws.Cells.Item(551, 4).Value = 45218
$ws.Cells.Item(551, 5).Value = 13
$ws.Cells.Item(551, 6).Value = 100112039
$ws.Cells.Item(551, 7).Value = "Ciboulette"
$ws.Cells.Item(551, 8).Value = "Sin especificar"
$ws.Cells.Item(551, 9).Value = "Primera"
$ws.Cells.Item(551, 10).Value = 250
$ws.Cells.Item(551, 11).Value = 1200
$ws.Cells.Item(551, 12).Value = 1500
$ws.Cells.Item(551, 13).Value = 1350
$ws.Cells.Item(551, 14).Value = "`$/docena de atados"
$ws.Cells.Item(551, 15).Value = "Región Metropolitana"
$ws.Cells.Item(551, 16).Value = 450
$ws.Cells.Item(551, 17).Value = 3
$ws.Cells.Item(551, 18).Value = "Hortaliza"
